$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 30, shifting existing rows 30-97 down to 31-98.
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with its data.
$ws.Cells.Item(30, 1).Value2 = 7
$ws.Cells.Item(30, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(30, 3).Value2 = "Ñuble"
$ws.Cells.Item(30, 4).Value2 = 45114
$ws.Cells.Item(30, 5).Value2 = 16
$ws.Cells.Item(30, 6).Value2 = 100112044
$ws.Cells.Item(30, 7).Value2 = "Perejil"
$ws.Cells.Item(30, 8).Value2 = "Sin especificar"
$ws.Cells.Item(30, 9).Value2 = "Primera"
$ws.Cells.Item(30, 10).Value2 = 100
$ws.Cells.Item(30, 11).Value2 = 1500
$ws.Cells.Item(30, 12).Value2 = 1500
$ws.Cells.Item(30, 13).Value2 = 1500
$ws.Cells.Item(30, 14).Value2 = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(30, 15).Value2 = "Región de Ñuble"
$ws.Cells.Item(30, 16).Value2 = 1500
$ws.Cells.Item(30, 17).Value2 = 1
$ws.Cells.Item(30, 18).Value2 = "Hortaliza"
